# #5: property boat&car done
# Complete the 汽車 (car) sheet: turn row 1 into a proper header row and
# extend rows 2-3 with the full set of metadata columns (H:N) that the
# other property sheets (land / building) already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Header row (row 1) --------------------------------------------------
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# ---- Row 2 (record #29) ---------------------------------------------------
$ws.Cells.Item(2, 1).Value = 29
$ws.Cells.Item(2, 2).Value = "NISSAN"
$ws.Cells.Item(2, 3).Value = 3500
$ws.Cells.Item(2, 4).Value = "邱議瑩"
$ws.Cells.Item(2, 5).Value = "93年06月18日"
$ws.Cells.Item(2, 6).Value = "買賣"
$ws.Cells.Item(2, 7).Value = "2000000(超過五年）"
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(2, 10).Value = "2011-11-25"
$ws.Cells.Item(2, 11).Value = "邱議瑩"
$ws.Cells.Item(2, 12).Value = 913
$ws.Cells.Item(2, 13).Value = "tmpab161"
$ws.Cells.Item(2, 14).Value = 29

# ---- Row 3 (record #30) ---------------------------------------------------
$ws.Cells.Item(3, 1).Value = 30
$ws.Cells.Item(3, 2).Value = "NISSAN"
$ws.Cells.Item(3, 3).Value = 3000
$ws.Cells.Item(3, 4).Value = "邱議瑩"
$ws.Cells.Item(3, 5).Value = "91年08月27日"
$ws.Cells.Item(3, 6).Value = "買賣"
$ws.Cells.Item(3, 7).Value = "(超過五年）"
$ws.Cells.Item(3, 8).Value = "land"
$ws.Cells.Item(3, 9).Value = "normal"
$ws.Cells.Item(3, 10).Value = "2011-11-25"
$ws.Cells.Item(3, 11).Value = "邱議瑩"
$ws.Cells.Item(3, 12).Value = 913
$ws.Cells.Item(3, 13).Value = "tmpab161"
$ws.Cells.Item(3, 14).Value = 30
